# "computed for covered? column and new car id"
#
# Adds the two derived columns to the car-inventory sheet:
#   M: "Covered?"  -> Y/N flag for whether Miles (H) are within Warantee Miles (L)
#   N: "New Car ID" -> a rebuilt car id from Make/Year/Model/Color/original id

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet already has the "Covered?" / "New Car ID" headers in M1/N1; only the
# new column (N) picked up an explicit custom width in the edit (M stayed default).
# ColumnWidth is in "characters" and gets re-quantized to pixels on save, so feed it
# a value that round-trips to the width closest to the original 17.26953125.
$ws.Columns.Item(14).ColumnWidth = 16.5

for ($r = 2; $r -le 53; $r++) {
    $ws.Range("M$r").Formula = "=IF(H$r<=L$r,""Y"", ""N"")"
    $ws.Range("N$r").Formula = "=CONCATENATE(B$r,F$r,D$r,UPPER(LEFT(J$r,3)),RIGHT(A$r,3))"
}

# Leave the selection where the author last left it when saving.
$ws.Range("N53").Select() | Out-Null
